# Scheduled-runner refresh of leve-profit price data across the ALC/ARM/
# BSM/CRP/CUL/GSM/LTW/WVR sheets (currentAveragePrice* / LevePrice* /
# LeveProfit* columns H-N). Values mirror a fresh market-board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4643.8184
$ws.Range("I74").Value = 3016.4
$ws.Range("K74").Value = 3016.4
$ws.Range("M74").Value = -2080.4
$ws.Range("H76").Value = 10899.389
$ws.Range("I76").Value = 12744.546
$ws.Range("J76").Value = 7999.857
$ws.Range("K76").Value = 12744.546
$ws.Range("L76").Value = 7999.857
$ws.Range("M76").Value = -12429.546
$ws.Range("N76").Value = -8629.857
$ws.Range("H77").Value = 4643.8184
$ws.Range("I77").Value = 3016.4
$ws.Range("K77").Value = 15082
$ws.Range("M77").Value = -10402
$ws.Range("H79").Value = 10899.389
$ws.Range("I79").Value = 12744.546
$ws.Range("J79").Value = 7999.857
$ws.Range("K79").Value = 12744.546
$ws.Range("L79").Value = 7999.857
$ws.Range("M79").Value = -11652.546
$ws.Range("N79").Value = -10183.857
$ws.Range("H100").Value = 1835
$ws.Range("I100").Value = 924.44446
$ws.Range("K100").Value = 924.44446
$ws.Range("M100").Value = -383.44446
$ws.Range("H113").Value = 9166.666999999999
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = 754
$ws.Range("H116").Value = 4501
$ws.Range("I116").Value = 4168.3335
$ws.Range("K116").Value = 4168.3335
$ws.Range("M116").Value = -726.3334999999997
$ws.Range("H132").Value = 2210.7551
$ws.Range("I132").Value = 1964
$ws.Range("K132").Value = 5892
$ws.Range("M132").Value = -3362
$ws.Range("H138").Value = 10102855
$ws.Range("I138").Value = 1241.5714
$ws.Range("J138").Value = 17546148
$ws.Range("K138").Value = 3724.7142
$ws.Range("L138").Value = 52638444
$ws.Range("M138").Value = 1415.2858
$ws.Range("N138").Value = -52648724
$ws.Range("H141").Value = 2112
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16399346
$ws.Range("I32").Value = 19234262
$ws.Range("J32").Value = 19833.334
$ws.Range("K32").Value = 19234262
$ws.Range("L32").Value = 19833.334
$ws.Range("M32").Value = -19233975
$ws.Range("N32").Value = -20407.334
$ws.Range("H63").Value = 4285.048
$ws.Range("I63").Value = 1859.1
$ws.Range("J63").Value = 6490.4546
$ws.Range("K63").Value = 1859.1
$ws.Range("L63").Value = 6490.4546
$ws.Range("M63").Value = -1173.1
$ws.Range("N63").Value = -7862.4546
$ws.Range("H66").Value = 4285.048
$ws.Range("I66").Value = 1859.1
$ws.Range("J66").Value = 6490.4546
$ws.Range("K66").Value = 9295.5
$ws.Range("L66").Value = 32452.273
$ws.Range("M66").Value = -5863.5
$ws.Range("N66").Value = -39316.273
$ws.Range("H132").Value = 25002766
$ws.Range("I132").Value = 2838.743
$ws.Range("K132").Value = 8516.228999999999
$ws.Range("M132").Value = -5986.228999999999
$ws.Range("H140").Value = 61888.8
$ws.Range("J140").Value = 61888.8
$ws.Range("L140").Value = 61888.8
$ws.Range("N140").Value = -72248.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4409.3335
$ws.Range("I107").Value = 4329.6665
$ws.Range("J107").Value = 4462.4443
$ws.Range("K107").Value = 4329.6665
$ws.Range("L107").Value = 4462.4443
$ws.Range("M107").Value = -2409.6665
$ws.Range("N107").Value = -8302.444299999999
$ws.Range("H134").Value = 3112.1428
$ws.Range("J134").Value = 6349.5
$ws.Range("L134").Value = 19048.5
$ws.Range("N134").Value = -24118.5
$ws.Range("H140").Value = 170775.39
$ws.Range("J140").Value = 170775.39
$ws.Range("L140").Value = 170775.39
$ws.Range("N140").Value = -181135.39

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23259858
$ws.Range("I31").Value = 2748.0667
$ws.Range("K31").Value = 2748.0667
$ws.Range("M31").Value = -2453.0667
$ws.Range("H32").Value = 716.3333
$ws.Range("I32").Value = 716.3333
$ws.Range("K32").Value = 716.3333
$ws.Range("M32").Value = -400.3333
$ws.Range("H34").Value = 23259858
$ws.Range("I34").Value = 2748.0667
$ws.Range("K34").Value = 2748.0667
$ws.Range("M34").Value = -2546.0667
$ws.Range("H58").Value = 2735.7368
$ws.Range("I58").Value = 1921.5385
$ws.Range("K58").Value = 1921.5385
$ws.Range("M58").Value = -1718.5385
$ws.Range("H62").Value = 3897.7
$ws.Range("I62").Value = 3449.5
$ws.Range("J62").Value = 4009.75
$ws.Range("K62").Value = 3449.5
$ws.Range("L62").Value = 4009.75
$ws.Range("M62").Value = -2825.5
$ws.Range("N62").Value = -5257.75
$ws.Range("H65").Value = 3897.7
$ws.Range("I65").Value = 3449.5
$ws.Range("J65").Value = 4009.75
$ws.Range("K65").Value = 17247.5
$ws.Range("L65").Value = 20048.75
$ws.Range("M65").Value = -14127.5
$ws.Range("N65").Value = -26288.75
$ws.Range("H86").Value = 4901.4
$ws.Range("I86").Value = 4502.3335
$ws.Range("J86").Value = 5500
$ws.Range("K86").Value = 4502.3335
$ws.Range("L86").Value = 5500
$ws.Range("M86").Value = -3379.3335
$ws.Range("N86").Value = -7746
$ws.Range("H89").Value = 4901.4
$ws.Range("I89").Value = 4502.3335
$ws.Range("J89").Value = 5500
$ws.Range("K89").Value = 22511.6675
$ws.Range("L89").Value = 27500
$ws.Range("M89").Value = -16895.6675
$ws.Range("N89").Value = -38732
$ws.Range("H94").Value = 1815.0714
$ws.Range("J94").Value = 2015.3334
$ws.Range("L94").Value = 2015.3334
$ws.Range("N94").Value = -2917.3334
$ws.Range("H132").Value = 3209.6333
$ws.Range("I132").Value = 2754.6316
$ws.Range("J132").Value = 3995.5454
$ws.Range("K132").Value = 8263.8948
$ws.Range("L132").Value = 11986.6362
$ws.Range("M132").Value = -5733.8948
$ws.Range("N132").Value = -17046.6362
$ws.Range("H136").Value = 2735.7368
$ws.Range("I136").Value = 1921.5385
$ws.Range("K136").Value = 5764.6155
$ws.Range("M136").Value = -3214.6155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 30079152
$ws.Range("I4").Value = 39004450
$ws.Range("J4").Value = 18030000
$ws.Range("K4").Value = 117013350
$ws.Range("L4").Value = 54090000
$ws.Range("M4").Value = -117013238
$ws.Range("N4").Value = -54090224
$ws.Range("H92").Value = 424.75
$ws.Range("J92").Value = 399.75
$ws.Range("L92").Value = 1199.25
$ws.Range("N92").Value = -3695.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4999.6665
$ws.Range("I80").Value = 4999.6665
$ws.Range("K80").Value = 4999.6665
$ws.Range("M80").Value = -4001.6665
$ws.Range("H83").Value = 4999.6665
$ws.Range("I83").Value = 4999.6665
$ws.Range("K83").Value = 24998.3325
$ws.Range("M83").Value = -20006.3325

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3928.5
$ws.Range("I7").Value = 4053
$ws.Range("J7").Value = 3736.0908
$ws.Range("K7").Value = 4053
$ws.Range("L7").Value = 3736.0908
$ws.Range("M7").Value = -3941
$ws.Range("N7").Value = -3960.0908
$ws.Range("H16").Value = 836
$ws.Range("I16").Value = 836
$ws.Range("K16").Value = 836
$ws.Range("M16").Value = -666
$ws.Range("H40").Value = 4109.857
$ws.Range("I40").Value = 4109.857
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4109.857
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3973.857
$ws.Range("H46").Value = 1418.2174
$ws.Range("I46").Value = 681
$ws.Range("J46").Value = 3103.2856
$ws.Range("K46").Value = 681
$ws.Range("L46").Value = 3103.2856
$ws.Range("M46").Value = -493
$ws.Range("N46").Value = -3479.2856
$ws.Range("H61").Value = 5350.2
$ws.Range("I61").Value = 4500
$ws.Range("K61").Value = 4500
$ws.Range("M61").Value = -4298
$ws.Range("H113").Value = 5350.2
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330
$ws.Range("H126").Value = 3928.5
$ws.Range("I126").Value = 4053
$ws.Range("J126").Value = 3736.0908
$ws.Range("K126").Value = 12159
$ws.Range("L126").Value = 11208.2724
$ws.Range("M126").Value = -9689
$ws.Range("N126").Value = -16148.2724
$ws.Range("H132").Value = 90912040
$ws.Range("I132").Value = 2942.5334
$ws.Range("K132").Value = 8827.600199999999
$ws.Range("M132").Value = -6297.600199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 626240.4
$ws.Range("I100").Value = 715528.4
$ws.Range("J100").Value = 1224.5
$ws.Range("K100").Value = 1431056.8
$ws.Range("L100").Value = 2449
$ws.Range("M100").Value = -1430515.8
$ws.Range("N100").Value = -3531
$ws.Range("H122").Value = 50001916
$ws.Range("I122").Value = 62501610
$ws.Range("K122").Value = 187504830
$ws.Range("M122").Value = -187502380
$ws.Range("H132").Value = 4648.852
$ws.Range("I132").Value = 4731.5386
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 14194.6158
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -11664.6158
$ws.Range("N132").Value = -12557
$ws.Range("H136").Value = 2065.8
$ws.Range("I136").Value = 1842.6111
$ws.Range("J136").Value = 2400.5833
$ws.Range("K136").Value = 5527.8333
$ws.Range("L136").Value = 7201.749899999999
$ws.Range("M136").Value = -2977.8333
$ws.Range("N136").Value = -12301.7499
